$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jobs_summary")

# Row 390 contains the job title "hr" which is a duplicate/unwanted entry
# (superseded by "hr generalist" / "human resources" rows). Delete it so the
# remaining rows shift up, matching the removal of the "hr" util helper data.
$ws.Rows.Item(390).Delete()

# Leave the selection where the deleted row used to be (matches the
# post-edit cursor position recorded in the workbook).
$ws.Range("A390").Select()
